$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2 should match the corrected value (same as B2/C2)
$ws.Range("D2").Value = 44545

# Rows 3-11: the "A" index column (10, 100, 11, 12, 13, 14, 15, 16, 17) was
# a spurious/missing-condition artifact and must be cleared, and the
# corresponding B/C/D measurement cells reset to 0.
$ws.Range("A3:A11").ClearContents()
$ws.Range("B3:D11").Value = 0
